# Auto-generated Excel COM-interop script
# Applies the cell-level numeric corrections described by the commit diff
# to the "Sophia_Profits" workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 62
$ws.Range("H62").Value = 2874.5
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 2999.3333
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 2999.3333
$ws.Range("M62").Value = -1876
$ws.Range("N62").Value = -4247.3333

# Row 65
$ws.Range("H65").Value = 2874.5
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 2999.3333
$ws.Range("K65").Value = 12500
$ws.Range("L65").Value = 14996.6665
$ws.Range("M65").Value = -9380
$ws.Range("N65").Value = -21236.6665

# Row 80
$ws.Range("H80").Value = 4549.5
$ws.Range("I80").Value = 3056.5715
$ws.Range("J80").Value = 15000
$ws.Range("K80").Value = 9169.7145
$ws.Range("L80").Value = 45000
$ws.Range("M80").Value = -8171.7145
$ws.Range("N80").Value = -46996

# Row 83
$ws.Range("H83").Value = 4549.5
$ws.Range("I83").Value = 3056.5715
$ws.Range("J83").Value = 15000
$ws.Range("K83").Value = 27509.1435
$ws.Range("L83").Value = 135000
$ws.Range("M83").Value = -22517.1435
$ws.Range("N83").Value = -144984

# Row 103
$ws.Range("H103").Value = 1915.5555
$ws.Range("I103").Value = 3000
$ws.Range("J103").Value = 1780
$ws.Range("K103").Value = 9000
$ws.Range("L103").Value = 5340
$ws.Range("M103").Value = -8414
$ws.Range("N103").Value = -6512

# Row 132
$ws.Range("H132").Value = 2301
$ws.Range("I132").Value = 1010.5455
$ws.Range("K132").Value = 3031.6365
$ws.Range("M132").Value = -501.6364999999996

# Row 137
$ws.Range("H137").Value = 3409.923
$ws.Range("J137").Value = 3537.9
$ws.Range("L137").Value = 10613.7
$ws.Range("N137").Value = -15713.7

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 74
$ws.Range("H74").Value = 6802.4
$ws.Range("I74").Value = 7878
$ws.Range("K74").Value = 7878
$ws.Range("M74").Value = -7004

# Row 77
$ws.Range("H77").Value = 6802.4
$ws.Range("I77").Value = 7878
$ws.Range("K77").Value = 39390
$ws.Range("M77").Value = -35022

# Row 122
$ws.Range("H122").Value = 2523
$ws.Range("I122").Value = 2045.4546
$ws.Range("J122").Value = 3573.6
$ws.Range("K122").Value = 6136.3638
$ws.Range("L122").Value = 10720.8
$ws.Range("M122").Value = -3686.3638
$ws.Range("N122").Value = -15620.8

# Row 132
$ws.Range("H132").Value = 3345.7693
$ws.Range("I132").Value = 2699.8
$ws.Range("K132").Value = 8099.400000000001
$ws.Range("M132").Value = -5569.400000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 20
$ws.Range("H20").Value = 3500
$ws.Range("I20").Value = 3500
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 3500
$ws.Range("L20").ClearContents()
$ws.Range("M20").Value = -3253
$ws.Range("N20").Value = 0

# Row 105
$ws.Range("H105").Value = 10187.875
$ws.Range("I105").Value = 11143.286
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 11143.286
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = -9396.286
$ws.Range("N105").Value = -6994

# Row 107
$ws.Range("H107").Value = 1165.7142
$ws.Range("I107").Value = 1015.25
$ws.Range("J107").Value = 1366.3334
$ws.Range("K107").Value = 1015.25
$ws.Range("L107").Value = 1366.3334
$ws.Range("M107").Value = 904.75
$ws.Range("N107").Value = -5206.3334

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 2439.2666
$ws.Range("J31").Value = 3379.2
$ws.Range("L31").Value = 3379.2
$ws.Range("N31").Value = -3969.2

# Row 34
$ws.Range("H34").Value = 2439.2666
$ws.Range("J34").Value = 3379.2
$ws.Range("L34").Value = 3379.2
$ws.Range("N34").Value = -3783.2

# Row 99
$ws.Range("H99").Value = 2525.9
$ws.Range("I99").Value = 2445
$ws.Range("J99").Value = 2849.5
$ws.Range("K99").Value = 2445
$ws.Range("L99").Value = 2849.5
$ws.Range("M99").Value = -947
$ws.Range("N99").Value = -5845.5

# Row 107
$ws.Range("H107").Value = 60043.47
$ws.Range("I107").Value = 84436.664
$ws.Range("K107").Value = 84436.664
$ws.Range("M107").Value = -82516.664

# Row 122
$ws.Range("H122").Value = 1081.2727
$ws.Range("I122").Value = 1296.5714
$ws.Range("J122").Value = 704.5
$ws.Range("K122").Value = 3889.7142
$ws.Range("L122").Value = 2113.5
$ws.Range("M122").Value = -1439.7142
$ws.Range("N122").Value = -7013.5

# Row 126
$ws.Range("H126").Value = 2525.9
$ws.Range("I126").Value = 2445
$ws.Range("J126").Value = 2849.5
$ws.Range("K126").Value = 7335
$ws.Range("L126").Value = 8548.5
$ws.Range("M126").Value = -4865
$ws.Range("N126").Value = -13488.5

# Row 132
$ws.Range("H132").Value = 104551
$ws.Range("J132").Value = 5358.7144
$ws.Range("L132").Value = 16076.1432
$ws.Range("N132").Value = -21136.1432

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

# Row 14
$ws.Range("H14").Value = 345
$ws.Range("I14").Value = 345
$ws.Range("K14").Value = 1035
$ws.Range("M14").Value = -862

# Row 52
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").ClearContents()
$ws.Range("N52").Value = 0

# Row 64
$ws.Range("H64").Value = 1000
$ws.Range("I64").Value = 1000
$ws.Range("K64").Value = 3000
$ws.Range("M64").Value = -2730

# Row 67
$ws.Range("H67").Value = 1000
$ws.Range("I67").Value = 1000
$ws.Range("K67").Value = 3000
$ws.Range("M67").Value = -2064

# Row 120
$ws.Range("H120").Value = 3500
$ws.Range("I120").Value = 3500
$ws.Range("K120").Value = 10500
$ws.Range("M120").Value = -5662

# Row 122
$ws.Range("H122").Value = 1569.875
$ws.Range("I122").Value = 1620.3334
$ws.Range("J122").Value = 1539.6
$ws.Range("K122").Value = 14583.0006
$ws.Range("L122").Value = 13856.4
$ws.Range("M122").Value = -12133.0006
$ws.Range("N122").Value = -18756.4

# Row 124
$ws.Range("H124").Value = 4583.5557
$ws.Range("I124").Value = 4030
$ws.Range("J124").Value = 4652.75
$ws.Range("K124").Value = 12090
$ws.Range("L124").Value = 13958.25
$ws.Range("M124").Value = -7180
$ws.Range("N124").Value = -23778.25

# Row 131
$ws.Range("H131").Value = 1672.3334
$ws.Range("I131").Value = 1641.1428
$ws.Range("J131").Value = 1699.625
$ws.Range("K131").Value = 4923.428400000001
$ws.Range("L131").Value = 5098.875
$ws.Range("M131").Value = 116.5715999999993
$ws.Range("N131").Value = -15178.875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 80
$ws.Range("H80").Value = 5800.375
$ws.Range("I80").Value = 2480.8
$ws.Range("K80").Value = 2480.8
$ws.Range("M80").Value = -1482.8

# Row 83
$ws.Range("H83").Value = 5800.375
$ws.Range("I83").Value = 2480.8
$ws.Range("K83").Value = 12404
$ws.Range("M83").Value = -7412

# Row 132
$ws.Range("H132").Value = 4749.7
$ws.Range("I132").Value = 3999.8333
$ws.Range("K132").Value = 11999.4999
$ws.Range("M132").Value = -9469.499899999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 10158.923
$ws.Range("I22").Value = 11595.143
$ws.Range("J22").Value = 8483.333000000001
$ws.Range("K22").Value = 11595.143
$ws.Range("L22").Value = 8483.333000000001
$ws.Range("M22").Value = -11300.143
$ws.Range("N22").Value = -9073.333000000001

# Row 27
$ws.Range("H27").Value = 10158.923
$ws.Range("I27").Value = 11595.143
$ws.Range("J27").Value = 8483.333000000001
$ws.Range("K27").Value = 11595.143
$ws.Range("L27").Value = 8483.333000000001
$ws.Range("M27").Value = -11488.143
$ws.Range("N27").Value = -8697.333000000001

# Row 55
$ws.Range("H55").Value = 552.5
$ws.Range("I55").Value = 631
$ws.Range("K55").Value = 631
$ws.Range("M55").Value = -458

# Row 122
$ws.Range("H122").Value = 7000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 21000
$ws.Range("N122").Value = -25900

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 96
$ws.Range("H96").Value = 3750
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

# Row 113
$ws.Range("H113").Value = 747.94116
$ws.Range("I113").Value = 740.1667
$ws.Range("J113").Value = 752.1818
$ws.Range("K113").Value = 2220.5001
$ws.Range("L113").Value = 2256.5454
$ws.Range("M113").Value = -50.5001000000002
$ws.Range("N113").Value = -6596.5454

# Row 126
$ws.Range("H126").Value = 4275
$ws.Range("I126").Value = 3742.8572
$ws.Range("K126").Value = 11228.5716
$ws.Range("M126").Value = -8758.571599999999

# Row 132
$ws.Range("H132").Value = 4419.6
$ws.Range("J132").Value = 4699.3335
$ws.Range("L132").Value = 14098.0005
$ws.Range("N132").Value = -19158.0005

